# Calibration of the needle: sort the curvature data rows (2-12) by the
# "time (s)" column (column A) in ascending order, keeping each row's
# A:D values together. Rows 6 and 9 already happen to be in the correct
# sorted position, so only the remaining rows are rewritten to avoid
# needlessly touching (and reformatting) cells that do not change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values (row number -> A,B,C,D), already sorted ascending by time.
# Values containing scientific notation are parsed via [double] cast from
# string literals because the script engine's numeric literal parser does
# not accept exponent suffixes (e.g. 1e-05) directly.
$rows = @{
    2  = @(56898.967185, [double]"-7.402089929e-06",      [double]"-1.1925646062e-05",      [double]"-2.5329835517e-05")
    3  = @(56916.699185, [double]"-4.7993912143e-05",     [double]"-0.00010193708406",      [double]"-9.9468749202e-05")
    4  = @(56927.699186, [double]"-5.6037888207e-05",     [double]"-0.00023315827516",      [double]"-0.00014170938367")
    5  = @(56938.567186, [double]"-7.5421499197e-05",     [double]"-0.00036139814919",      [double]"-0.00019064895824")
    7  = @(56960.367188, [double]"-0.0002014021",         [double]"-0.0006097965",          [double]"-0.0003018801")
    8  = @(56978.167188, [double]"-0.0001424692",         [double]"-0.0004881551",          [double]"-0.0002451157")
    10 = @(57001.231189, [double]"-5.5720426976e-05",     [double]"-0.00022620906777",      [double]"-0.00014607513843")
    11 = @(57011.56719,  [double]"-2.4025626554e-05",     [double]"-7.608194195399999e-05", [double]"-0.00010670966426")
    12 = @(57021.63119,  [double]"-7.2268442087e-06",     [double]"-1.1698795606e-05",      [double]"-2.495058565e-05")
}

foreach ($r in $rows.Keys) {
    $entry = $rows[$r]
    $ws.Cells.Item($r, 1).Value = $entry[0]
    $ws.Cells.Item($r, 2).Value = $entry[1]
    $ws.Cells.Item($r, 3).Value = $entry[2]
    $ws.Cells.Item($r, 4).Value = $entry[3]
}
